$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 0.01451666666666667
$ws.Cells.Item(2, 8).Value = 0.04355
$ws.Cells.Item(2, 9).Value = 0.09504707612932513
$ws.Cells.Item(2, 10).Value = 0.09504707612932513
$ws.Cells.Item(2, 13).Value = 0.036961
$ws.Cells.Item(2, 14).Value = 0.110883
$ws.Cells.Item(2, 15).Value = 0.1786608532565087
$ws.Cells.Item(2, 16).Value = 0.1786608532565087
$ws.Cells.Item(2, 17).Value = 0.0005365505166666666
$ws.Cells.Item(2, 18).Value = 0.00482895465
$ws.Cells.Item(2, 19).Value = 0.01698119172080156
$ws.Cells.Item(2, 20).Value = 0.01698119172080156
$ws.Cells.Item(3, 7).Value = 0.01451666666666667
$ws.Cells.Item(3, 8).Value = 0.04355
$ws.Cells.Item(3, 9).Value = 0.09504707612932513
$ws.Cells.Item(3, 10).Value = 0.09504707612932513
$ws.Cells.Item(3, 14).Value = 0.384071
$ws.Cells.Item(3, 15).Value = 0.6188365445657183
$ws.Cells.Item(3, 16).Value = 0.6188365445657182
$ws.Cells.Item(3, 17).Value = 0.001858476894444444
$ws.Cells.Item(3, 18).Value = 0.01672629205
$ws.Cells.Item(3, 19).Value = 0.05881860416294633
$ws.Cells.Item(3, 20).Value = 0.05881860416294633
$ws.Cells.Item(4, 7).Value = 0.01451666666666667
$ws.Cells.Item(4, 8).Value = 0.04355
$ws.Cells.Item(4, 9).Value = 0.09504707612932513
$ws.Cells.Item(4, 10).Value = 0.09504707612932513
$ws.Cells.Item(4, 15).Value = 0.202502602177773
$ws.Cells.Item(4, 16).Value = 0.202502602177773
$ws.Cells.Item(4, 17).Value = 0.0006081515555555555
$ws.Cells.Item(4, 18).Value = 0.005473363999999999
$ws.Cells.Item(4, 19).Value = 0.01924728024557723
$ws.Cells.Item(4, 20).Value = 0.01924728024557723
$ws.Cells.Item(5, 9).Value = 0.5165388459909994
$ws.Cells.Item(5, 10).Value = 0.5165388459909994
$ws.Cells.Item(5, 13).Value = 0.036961
$ws.Cells.Item(5, 14).Value = 0.110883
$ws.Cells.Item(5, 15).Value = 0.1786608532565087
$ws.Cells.Item(5, 16).Value = 0.1786608532565087
$ws.Cells.Item(5, 17).Value = 0.002915914891666667
$ws.Cells.Item(5, 18).Value = 0.026243234025
$ws.Cells.Item(5, 19).Value = 0.09228527096488427
$ws.Cells.Item(5, 20).Value = 0.09228527096488427
$ws.Cells.Item(6, 9).Value = 0.5165388459909994
$ws.Cells.Item(6, 10).Value = 0.5165388459909994
$ws.Cells.Item(6, 14).Value = 0.384071
$ws.Cells.Item(6, 15).Value = 0.6188365445657183
$ws.Cells.Item(6, 16).Value = 0.6188365445657182
$ws.Cells.Item(6, 18).Value = 0.090900003925
$ws.Cells.Item(6, 19).Value = 0.3196531145870338
$ws.Cells.Item(6, 20).Value = 0.3196531145870338
$ws.Cells.Item(7, 9).Value = 0.5165388459909994
$ws.Cells.Item(7, 10).Value = 0.5165388459909994
$ws.Cells.Item(7, 15).Value = 0.202502602177773
$ws.Cells.Item(7, 16).Value = 0.202502602177773
$ws.Cells.Item(7, 19).Value = 0.1046004604390813
$ws.Cells.Item(7, 20).Value = 0.1046004604390813
$ws.Cells.Item(8, 9).Value = 0.3884140778796754
$ws.Cells.Item(8, 10).Value = 0.3884140778796754
$ws.Cells.Item(8, 13).Value = 0.036961
$ws.Cells.Item(8, 14).Value = 0.110883
$ws.Cells.Item(8, 15).Value = 0.1786608532565087
$ws.Cells.Item(8, 16).Value = 0.1786608532565087
$ws.Cells.Item(8, 17).Value = 0.002192637403
$ws.Cells.Item(8, 18).Value = 0.019733736627
$ws.Cells.Item(8, 19).Value = 0.06939439057082282
$ws.Cells.Item(8, 20).Value = 0.06939439057082282
$ws.Cells.Item(9, 9).Value = 0.3884140778796754
$ws.Cells.Item(9, 10).Value = 0.3884140778796754
$ws.Cells.Item(9, 14).Value = 0.384071
$ws.Cells.Item(9, 15).Value = 0.6188365445657183
$ws.Cells.Item(9, 16).Value = 0.6188365445657182
$ws.Cells.Item(9, 18).Value = 0.06835273179899999
$ws.Cells.Item(9, 19).Value = 0.2403648258157381
$ws.Cells.Item(9, 20).Value = 0.2403648258157381
$ws.Cells.Item(10, 9).Value = 0.3884140778796754
$ws.Cells.Item(10, 10).Value = 0.3884140778796754
$ws.Cells.Item(10, 15).Value = 0.202502602177773
$ws.Cells.Item(10, 16).Value = 0.202502602177773
$ws.Cells.Item(10, 19).Value = 0.07865486149311446
$ws.Cells.Item(10, 20).Value = 0.07865486149311446
